$wb = $excel.ActiveWorkbook

# Mapping of row -> [oldValue, newValue] for column F ("想去人数")
# applied identically to both the "展览" and "全部类型" worksheets.
$updates = @(
    @{Row = 3;  New = 495},
    @{Row = 4;  New = 1265},
    @{Row = 5;  New = 1098},
    @{Row = 6;  New = 14099},
    @{Row = 7;  New = 15662},
    @{Row = 11; New = 190},
    @{Row = 17; New = 30},
    @{Row = 19; New = 29},
    @{Row = 20; New = 1216},
    @{Row = 21; New = 128},
    @{Row = 23; New = 6143},
    @{Row = 24; New = 961},
    @{Row = 25; New = 1093},
    @{Row = 26; New = 5555},
    @{Row = 29; New = 118},
    @{Row = 30; New = 4530}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.New
}

# "全部类型" sheet mirrors the same data, but has one extra row inserted
# before row 23, so every row from 23 onward is shifted down by one.
$updates4 = @(
    @{Row = 3;  New = 495},
    @{Row = 4;  New = 1265},
    @{Row = 5;  New = 1098},
    @{Row = 6;  New = 14099},
    @{Row = 7;  New = 15662},
    @{Row = 11; New = 190},
    @{Row = 17; New = 30},
    @{Row = 19; New = 29},
    @{Row = 20; New = 1216},
    @{Row = 21; New = 128},
    @{Row = 24; New = 6143},
    @{Row = 25; New = 961},
    @{Row = 26; New = 1093},
    @{Row = 27; New = 5555},
    @{Row = 30; New = 118},
    @{Row = 31; New = 4530}
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.New
}
